$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 38; this shifts the existing rows 38-68 down to 39-69.
$ws.Rows("38").Insert()

# Populate the newly inserted row 38 with the new weekly record.
$ws.Range("A38").Value = 7
$ws.Range("B38").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C38").Value = "Ñuble"
$ws.Range("D38").Value = 45040
$ws.Range("D38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E38").Value = 16
$ws.Range("F38").Value = 100112001
$ws.Range("G38").Value = "Berenjena"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 50
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("M38").Value = 10000
$ws.Range("N38").Value = "$/caja 60 unidades"
$ws.Range("O38").Value = "Región de Arica y Parinacota"
$ws.Range("P38").Value = 167
$ws.Range("Q38").Value = 60
$ws.Range("R38").Value = "Hortaliza"
